$d = $word.ActiveDocument
$last = $d.Paragraphs.Last
$r = $last.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Daily</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Scrum - Día 3</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Facilitador: Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Alán Osmar Peña Polo (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué hice ayer?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Completé el componente de visualización de resultados (T2.4) </w:t></w:r></w:p><w:p><w:r><w:t>Modal responsive implementado</w:t></w:r></w:p><w:p><w:r><w:t>Animaciones CSS funcionando correctamente</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Efecto de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>confetti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para victorias (¡se ve increíble!)</w:t></w:r></w:p><w:p><w:r><w:t>Diseño visualmente atractivo</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Colaboré en las pruebas de renderizado (T2.5)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Validamos que los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> se ven perfectamente en el LCD</w:t></w:r></w:p><w:p><w:r><w:t>Documentamos 8 casos de prueba visuales</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementar recepción de telemetría en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flask</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T3.6)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Crear </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endpoint</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>listener</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para recibir JSON del PIC</w:t></w:r></w:p><w:p><w:r><w:t>Validar estructura de telemetría recibida</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Preparar envío al </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>frontend</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">Comenzar implementación de comunicación </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>frontend-backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para telemetría (avance de T4.3)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Evaluar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>polling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> vs </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WebSocket</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">Probablemente implementar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>polling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (más simple)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t>Ninguno. El componente de visualización está listo para recibir datos reales del PIC.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Embebido)</w:t></w:r></w:p><w:p><w:r><w:t>¿Qué hice ayer?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Implementé inicialización del juego con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T2.1)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> se cargan correctamente de CGRAM</w:t></w:r></w:p><w:p><w:r><w:t>El personaje aparece en la primera columna</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Implementé el sistema de movimiento del personaje (T2.2) </w:t></w:r></w:p><w:p><w:r><w:t>2 botones configurados para 2 carriles (fila superior e inferior)</w:t></w:r></w:p><w:p><w:r><w:t>El movimiento es responsive y fluido</w:t></w:r></w:p><w:p><w:r><w:t>El renderizado funciona sin problemas</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Implementé generación y movimiento de obstáculos (T2.3) </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sistema de generación </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pseudo-aleatoria</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> funcionando</w:t></w:r></w:p><w:p><w:r><w:t>Buffer de 10 obstáculos activos implementado</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Desplazamiento horizontal funcional</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Colaboramos en pruebas de renderizado (T2.5) </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p><w:r><w:t>Implementar algoritmo de detección de colisiones (T3.1)</w:t></w:r></w:p><w:p><w:r><w:t>Comparar posiciones del personaje vs obstáculos</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Definir </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hitboxes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (probablemente 5x8 completo para precisión)</w:t></w:r></w:p><w:p><w:r><w:t>Manejar casos límite (obstáculo saliendo de pantalla, múltiples obstáculos)</w:t></w:r></w:p><w:p><w:r><w:t>Implementar evaluación de metas por obstáculos esquivados (T3.2)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Contador de obstáculos que pasan sin colisión</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Comparar con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cuando </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> == 0</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Implementar evaluación de metas por tiempo sobrevivido (T3.3)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Configurar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Timer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> con interrupciones para contar segundos</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Comparar con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cuando </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> == 1</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t>IMPEDIMENTO CRÍTICO: Durante las pruebas de ayer noche, el uso de memoria RAM alcanzó 352/368 bytes (95.6%). El buffer de 10 obstáculos consume demasiada RAM. Necesito optimizar urgentemente antes de continuar con colisiones.</w:t></w:r></w:p></w:body></w:wordDocument>'
$r.InsertXML($xml)
Write-Output "Inserted new Daily Scrum Dia 3 content"
